$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows (82-137) for the newly-added data sets
$ws.Cells.Item(82, 1).Value = "CH4"
$ws.Cells.Item(82, 2).Value = 139.30464912377872
$ws.Cells.Item(82, 3).Value = 921.6391379649823
$ws.Cells.Item(82, 4).Value = 2
$ws.Cells.Item(83, 1).Value = "CH8"
$ws.Cells.Item(83, 2).Value = 133.51868036160101
$ws.Cells.Item(83, 3).Value = 682.65556698579053
$ws.Cells.Item(83, 4).Value = 1
$ws.Cells.Item(84, 1).Value = "CH12"
$ws.Cells.Item(84, 2).Value = 92.875948209028977
$ws.Cells.Item(84, 3).Value = 714.06440672507654
$ws.Cells.Item(84, 4).Value = 1
$ws.Cells.Item(85, 1).Value = "ORIGINAL"
$ws.Cells.Item(85, 2).Value = 114.37112615658687
$ws.Cells.Item(85, 3).Value = 422.89784820263202
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 1).Value = "ORIGINAL"
$ws.Cells.Item(86, 2).Value = 85.312759784551773
$ws.Cells.Item(86, 3).Value = 457.77465798304632
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 1).Value = "CH4"
$ws.Cells.Item(87, 2).Value = 59.447127470603355
$ws.Cells.Item(87, 3).Value = 264.57460381434515
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 1).Value = "CH8"
$ws.Cells.Item(88, 2).Value = 134.35369909726657
$ws.Cells.Item(88, 3).Value = 331.53298466022198
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 1).Value = "CH12"
$ws.Cells.Item(89, 2).Value = 60.747052394426788
$ws.Cells.Item(89, 3).Value = 278.24108285170337
$ws.Cells.Item(89, 4).Value = 0
$ws.Cells.Item(90, 1).Value = "CH4"
$ws.Cells.Item(90, 2).Value = 17.734175828786995
$ws.Cells.Item(90, 3).Value = 15.76536618746244
$ws.Cells.Item(90, 4).Value = 3
$ws.Cells.Item(91, 1).Value = "CH8"
$ws.Cells.Item(91, 2).Value = 16.679549437302811
$ws.Cells.Item(91, 3).Value = 17.642024553739109
$ws.Cells.Item(91, 4).Value = 2
$ws.Cells.Item(92, 1).Value = "CH12"
$ws.Cells.Item(92, 2).Value = 19.761544594397911
$ws.Cells.Item(92, 3).Value = 21.199544173020584
$ws.Cells.Item(92, 4).Value = 1
$ws.Cells.Item(93, 1).Value = "ORIGINAL"
$ws.Cells.Item(93, 2).Value = 20.361760799701397
$ws.Cells.Item(93, 3).Value = 17.875503466679501
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 1).Value = "ORIGINAL"
$ws.Cells.Item(94, 2).Value = 21.0097226362962
$ws.Cells.Item(94, 3).Value = 15.958862231327938
$ws.Cells.Item(94, 4).Value = 0
$ws.Cells.Item(95, 1).Value = "CH4"
$ws.Cells.Item(95, 2).Value = 42.094740794255181
$ws.Cells.Item(95, 3).Value = 15.279813913198618
$ws.Cells.Item(95, 4).Value = 0
$ws.Cells.Item(96, 1).Value = "CH8"
$ws.Cells.Item(96, 2).Value = 19.872208081758938
$ws.Cells.Item(96, 3).Value = 13.41997109926664
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(97, 1).Value = "CH12"
$ws.Cells.Item(97, 2).Value = 32.27365537790152
$ws.Cells.Item(97, 3).Value = 22.17827320098877
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 1).Value = "CH4"
$ws.Cells.Item(98, 2).Value = 10.984649364764874
$ws.Cells.Item(98, 3).Value = 7.3801723844729938
$ws.Cells.Item(98, 4).Value = 2
$ws.Cells.Item(99, 1).Value = "CH8"
$ws.Cells.Item(99, 2).Value = 43.437601749713608
$ws.Cells.Item(99, 3).Value = 17.10391715856699
$ws.Cells.Item(99, 4).Value = 2
$ws.Cells.Item(100, 1).Value = "CH12"
$ws.Cells.Item(100, 2).Value = 15.946001254595243
$ws.Cells.Item(100, 3).Value = 10.220108740604841
$ws.Cells.Item(100, 4).Value = 2
$ws.Cells.Item(101, 1).Value = "ORIGINAL"
$ws.Cells.Item(101, 2).Value = 12.259904678051289
$ws.Cells.Item(101, 3).Value = 10.030653958137218
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(102, 1).Value = "ORIGINAL"
$ws.Cells.Item(102, 2).Value = 43.637437306917633
$ws.Cells.Item(102, 3).Value = 23.242417107407864
$ws.Cells.Item(102, 4).Value = 0
$ws.Cells.Item(103, 1).Value = "CH4"
$ws.Cells.Item(103, 2).Value = 50.297773810533378
$ws.Cells.Item(103, 3).Value = 10.180087707936764
$ws.Cells.Item(103, 4).Value = 0
$ws.Cells.Item(104, 1).Value = "CH8"
$ws.Cells.Item(104, 2).Value = 11.398828423940218
$ws.Cells.Item(104, 3).Value = 9.0864980197869816
$ws.Cells.Item(104, 4).Value = 0
$ws.Cells.Item(105, 1).Value = "CH12"
$ws.Cells.Item(105, 2).Value = 34.465603702343429
$ws.Cells.Item(105, 3).Value = 9.6959857677037906
$ws.Cells.Item(105, 4).Value = 0
$ws.Cells.Item(106, 1).Value = "CH4"
$ws.Cells.Item(106, 2).Value = 37.577071432883926
$ws.Cells.Item(106, 3).Value = 10.336478711320805
$ws.Cells.Item(106, 4).Value = 3
$ws.Cells.Item(107, 1).Value = "CH8"
$ws.Cells.Item(107, 2).Value = 20.580863154851475
$ws.Cells.Item(107, 3).Value = 7.2496556421885128
$ws.Cells.Item(107, 4).Value = 1
$ws.Cells.Item(108, 1).Value = "CH12"
$ws.Cells.Item(108, 2).Value = 45.140778720378876
$ws.Cells.Item(108, 3).Value = 6.7889590011193199
$ws.Cells.Item(108, 4).Value = 1
$ws.Cells.Item(109, 1).Value = "ORIGINAL"
$ws.Cells.Item(109, 2).Value = 35.592536733700676
$ws.Cells.Item(109, 3).Value = 6.8021187243553305
$ws.Cells.Item(109, 4).Value = 0
$ws.Cells.Item(110, 1).Value = "ORIGINAL"
$ws.Cells.Item(110, 2).Value = 46.739307637398056
$ws.Cells.Item(110, 3).Value = 10.987723417007006
$ws.Cells.Item(110, 4).Value = 0
$ws.Cells.Item(111, 1).Value = "CH4"
$ws.Cells.Item(111, 2).Value = 35.884697350171898
$ws.Cells.Item(111, 3).Value = 11.033881172537804
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(112, 1).Value = "CH8"
$ws.Cells.Item(112, 2).Value = 26.148272943038208
$ws.Cells.Item(112, 3).Value = 7.5593141718552666
$ws.Cells.Item(112, 4).Value = 0
$ws.Cells.Item(113, 1).Value = "CH12"
$ws.Cells.Item(113, 2).Value = 26.662141859531403
$ws.Cells.Item(113, 3).Value = 8.6836000428749962
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(114, 1).Value = "CH4"
$ws.Cells.Item(114, 2).Value = 66.342754134765045
$ws.Cells.Item(114, 3).Value = 13.0430716918065
$ws.Cells.Item(114, 4).Value = 1
$ws.Cells.Item(115, 1).Value = "CH8"
$ws.Cells.Item(115, 2).Value = 67.357748590982879
$ws.Cells.Item(115, 3).Value = 10.095222266820761
$ws.Cells.Item(115, 4).Value = 1
$ws.Cells.Item(116, 1).Value = "CH12"
$ws.Cells.Item(116, 2).Value = 61.879314853594856
$ws.Cells.Item(116, 3).Value = 8.1855825999608403
$ws.Cells.Item(116, 4).Value = 1
$ws.Cells.Item(117, 1).Value = "ORIGINAL"
$ws.Cells.Item(117, 2).Value = 46.610791355371475
$ws.Cells.Item(117, 3).Value = 11.350654348731041
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(118, 1).Value = "ORIGINAL"
$ws.Cells.Item(118, 2).Value = 51.166381372855263
$ws.Cells.Item(118, 3).Value = 8.4160332301488285
$ws.Cells.Item(118, 4).Value = 0
$ws.Cells.Item(119, 1).Value = "CH4"
$ws.Cells.Item(119, 2).Value = 62.964724467350884
$ws.Cells.Item(119, 3).Value = 11.448315840501051
$ws.Cells.Item(119, 4).Value = 0
$ws.Cells.Item(120, 1).Value = "CH8"
$ws.Cells.Item(120, 2).Value = 52.84523785802034
$ws.Cells.Item(120, 3).Value = 9.9696730535763969
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(121, 1).Value = "CH12"
$ws.Cells.Item(121, 2).Value = 77.908704138719116
$ws.Cells.Item(121, 3).Value = 9.3947173609183388
$ws.Cells.Item(121, 4).Value = 0
$ws.Cells.Item(122, 1).Value = "CH4"
$ws.Cells.Item(122, 2).Value = 32.858189431520607
$ws.Cells.Item(122, 3).Value = 5.0025127931283073
$ws.Cells.Item(122, 4).Value = 3
$ws.Cells.Item(123, 1).Value = "CH8"
$ws.Cells.Item(123, 2).Value = 51.641768840643074
$ws.Cells.Item(123, 3).Value = 4.5084924984436769
$ws.Cells.Item(123, 4).Value = 1
$ws.Cells.Item(124, 1).Value = "CH12"
$ws.Cells.Item(124, 2).Value = 50.78301159005899
$ws.Cells.Item(124, 3).Value = 9.5494402394844933
$ws.Cells.Item(124, 4).Value = 1
$ws.Cells.Item(125, 1).Value = "ORIGINAL"
$ws.Cells.Item(125, 2).Value = 28.867430583788799
$ws.Cells.Item(125, 3).Value = 4.8104356378316879
$ws.Cells.Item(125, 4).Value = 0
$ws.Cells.Item(126, 1).Value = "ORIGINAL"
$ws.Cells.Item(126, 2).Value = 21.107921126943367
$ws.Cells.Item(126, 3).Value = 3.8491739447300253
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(127, 1).Value = "CH4"
$ws.Cells.Item(127, 2).Value = 13.486887950163622
$ws.Cells.Item(127, 3).Value = 2.8865842211704988
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(128, 1).Value = "CH8"
$ws.Cells.Item(128, 2).Value = 18.927964623157795
$ws.Cells.Item(128, 3).Value = 4.0615977335434694
$ws.Cells.Item(128, 4).Value = 0
$ws.Cells.Item(129, 1).Value = "CH12"
$ws.Cells.Item(129, 2).Value = 16.830797405197071
$ws.Cells.Item(129, 3).Value = 3.2112299971855602
$ws.Cells.Item(129, 4).Value = 0
$ws.Cells.Item(130, 1).Value = "CH4"
$ws.Cells.Item(130, 2).Value = 19.396602043738731
$ws.Cells.Item(130, 3).Value = 11.061275019095493
$ws.Cells.Item(130, 4).Value = 2
$ws.Cells.Item(131, 1).Value = "CH8"
$ws.Cells.Item(131, 2).Value = 36.426098080781792
$ws.Cells.Item(131, 3).Value = 9.3817504094197197
$ws.Cells.Item(131, 4).Value = 1
$ws.Cells.Item(132, 1).Value = "CH12"
$ws.Cells.Item(132, 2).Value = 40.917258088405319
$ws.Cells.Item(132, 3).Value = 12.612592224891369
$ws.Cells.Item(132, 4).Value = 1
$ws.Cells.Item(133, 1).Value = "ORIGINAL"
$ws.Cells.Item(133, 2).Value = 28.658360536281879
$ws.Cells.Item(133, 3).Value = 9.9017491432336655
$ws.Cells.Item(133, 4).Value = 0
$ws.Cells.Item(134, 1).Value = "ORIGINAL"
$ws.Cells.Item(134, 2).Value = 34.38826078634996
$ws.Cells.Item(134, 3).Value = 13.508753666510948
$ws.Cells.Item(134, 4).Value = 0
$ws.Cells.Item(135, 1).Value = "CH4"
$ws.Cells.Item(135, 2).Value = 29.271822892702541
$ws.Cells.Item(135, 3).Value = 11.188973656067482
$ws.Cells.Item(135, 4).Value = 0
$ws.Cells.Item(136, 1).Value = "CH8"
$ws.Cells.Item(136, 2).Value = 29.971125052525448
$ws.Cells.Item(136, 3).Value = 11.510829017712521
$ws.Cells.Item(136, 4).Value = 0
$ws.Cells.Item(137, 1).Value = "CH12"
$ws.Cells.Item(137, 2).Value = 41.808013365818901
$ws.Cells.Item(137, 3).Value = 11.297479125169607
$ws.Cells.Item(137, 4).Value = 0

# Update the view selection to match the end of the newly added data
$ws.Range("G102").Select()
